$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 7.312497333333334
$ws.Range("H2").Value = 21.937492
$ws.Range("I2").Value = 0.05970572560549242
$ws.Range("J2").Value = 0.05970572560549242
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.752878
$ws.Range("N2").Value = 5.258634
$ws.Range("O2").Value = 0.1377607590022273
$ws.Range("P2").Value = 0.1377607590022273
$ws.Range("Q2").Value = 12.81791570065867
$ws.Range("R2").Value = 115.361241305928
$ws.Range("S2").Value = 0.008225106076191353
$ws.Range("T2").Value = 0.008225106076191353
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 7.312497333333334
$ws.Range("H3").Value = 21.937492
$ws.Range("I3").Value = 0.05970572560549242
$ws.Range("J3").Value = 0.05970572560549242
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.076282333333333
$ws.Range("N3").Value = 9.228847
$ws.Range("O3").Value = 0.2417686736584878
$ws.Range("P3").Value = 0.2417686736584878
$ws.Range("Q3").Value = 22.49530635908044
$ws.Range("R3").Value = 202.457757231724
$ws.Range("S3").Value = 0.01443497408945752
$ws.Range("T3").Value = 0.01443497408945752
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 7.312497333333334
$ws.Range("H4").Value = 21.937492
$ws.Range("I4").Value = 0.05970572560549242
$ws.Range("J4").Value = 0.05970572560549242
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.303088666666667
$ws.Range("N4").Value = 21.909266
$ws.Range("O4").Value = 0.5739583917309499
$ws.Range("P4").Value = 0.5739583917309499
$ws.Range("Q4").Value = 53.40381640009689
$ws.Range("R4").Value = 480.634347600872
$ws.Range("S4").Value = 0.03426860224565782
$ws.Range("T4").Value = 0.03426860224565782
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 7.312497333333334
$ws.Range("H5").Value = 21.937492
$ws.Range("I5").Value = 0.05970572560549242
$ws.Range("J5").Value = 0.05970572560549242
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.5918243333333334
$ws.Range("N5").Value = 1.775473
$ws.Range("O5").Value = 0.04651217560833507
$ws.Range("P5").Value = 0.04651217560833507
$ws.Range("Q5").Value = 4.327713859301778
$ws.Range("R5").Value = 38.94942473371601
$ws.Range("S5").Value = 0.002777043194185731
$ws.Range("T5").Value = 0.002777043194185731
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 29.68221266666667
$ws.Range("H6").Value = 89.046638
$ws.Range("I6").Value = 0.242351957758873
$ws.Range("J6").Value = 0.242351957758873
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.752878
$ws.Range("N6").Value = 5.258634
$ws.Range("O6").Value = 0.1377607590022273
$ws.Range("P6").Value = 0.1377607590022273
$ws.Range("Q6").Value = 52.02929757472133
$ws.Range("R6").Value = 468.263678172492
$ws.Range("S6").Value = 0.03338658964653807
$ws.Range("T6").Value = 0.03338658964653807
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 29.68221266666667
$ws.Range("H7").Value = 89.046638
$ws.Range("I7").Value = 0.242351957758873
$ws.Range("J7").Value = 0.242351957758873
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.076282333333333
$ws.Range("N7").Value = 9.228847
$ws.Range("O7").Value = 0.2417686736584878
$ws.Range("P7").Value = 0.2417686736584878
$ws.Range("Q7").Value = 91.31086644070956
$ws.Range("R7").Value = 821.797797966386
$ws.Range("S7").Value = 0.05859311138590058
$ws.Range("T7").Value = 0.05859311138590058
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 29.68221266666667
$ws.Range("H8").Value = 89.046638
$ws.Range("I8").Value = 0.242351957758873
$ws.Range("J8").Value = 0.242351957758873
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 7.303088666666667
$ws.Range("N8").Value = 21.909266
$ws.Range("O8").Value = 0.5739583917309499
$ws.Range("P8").Value = 0.5739583917309499
$ws.Range("Q8").Value = 216.7718309275231
$ws.Range("R8").Value = 1950.946478347708
$ws.Range("S8").Value = 0.1390999399081299
$ws.Range("T8").Value = 0.1390999399081298
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 29.68221266666667
$ws.Range("H9").Value = 89.046638
$ws.Range("I9").Value = 0.242351957758873
$ws.Range("J9").Value = 0.242351957758873
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.5918243333333334
$ws.Range("N9").Value = 1.775473
$ws.Range("O9").Value = 0.04651217560833507
$ws.Range("P9").Value = 0.04651217560833507
$ws.Range("Q9").Value = 17.56665572330823
$ws.Range("R9").Value = 158.099901509774
$ws.Range("S9").Value = 0.0112723168183045
$ws.Range("T9").Value = 0.0112723168183045
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 11.06470466666667
$ws.Range("H10").Value = 33.194114
$ws.Range("I10").Value = 0.09034208022509747
$ws.Range("J10").Value = 0.09034208022509749
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.752878
$ws.Range("N10").Value = 5.258634
$ws.Range("O10").Value = 0.1377607590022273
$ws.Range("P10").Value = 0.1377607590022273
$ws.Range("Q10").Value = 19.39507738669733
$ws.Range("R10").Value = 174.555696480276
$ws.Range("S10").Value = 0.01244559354164954
$ws.Range("T10").Value = 0.01244559354164954
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 11.06470466666667
$ws.Range("H11").Value = 33.194114
$ws.Range("I11").Value = 0.09034208022509747
$ws.Range("J11").Value = 0.09034208022509749
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 3.076282333333333
$ws.Range("N11").Value = 9.228847
$ws.Range("O11").Value = 0.2417686736584878
$ws.Range("P11").Value = 0.2417686736584878
$ws.Range("Q11").Value = 34.03815548961755
$ws.Range("R11").Value = 306.343399406558
$ws.Range("S11").Value = 0.02184188491157051
$ws.Range("T11").Value = 0.02184188491157052
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 11.06470466666667
$ws.Range("H12").Value = 33.194114
$ws.Range("I12").Value = 0.09034208022509747
$ws.Range("J12").Value = 0.09034208022509749
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 7.303088666666667
$ws.Range("N12").Value = 21.909266
$ws.Range("O12").Value = 0.5739583917309499
$ws.Range("P12").Value = 0.5739583917309499
$ws.Range("Q12").Value = 80.80651925114711
$ws.Range("R12").Value = 727.2586732603239
$ws.Range("S12").Value = 0.05185259507162539
$ws.Range("T12").Value = 0.0518525950716254
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 11.06470466666667
$ws.Range("H13").Value = 33.194114
$ws.Range("I13").Value = 0.09034208022509747
$ws.Range("J13").Value = 0.09034208022509749
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.5918243333333334
$ws.Range("N13").Value = 1.775473
$ws.Range("O13").Value = 0.04651217560833507
$ws.Range("P13").Value = 0.04651217560833507
$ws.Range("Q13").Value = 6.548361462880222
$ws.Range("R13").Value = 58.935253165922
$ws.Range("S13").Value = 0.004202006700252029
$ws.Range("T13").Value = 0.004202006700252029
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 74.41623166666666
$ws.Range("H14").Value = 223.248695
$ws.Range("I14").Value = 0.6076002364105371
$ws.Range("J14").Value = 0.6076002364105371
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.752878
$ws.Range("N14").Value = 5.258634
$ws.Range("O14").Value = 0.1377607590022273
$ws.Range("P14").Value = 0.1377607590022273
$ws.Range("Q14").Value = 130.4425753314033
$ws.Range("R14").Value = 1173.98317798263
$ws.Range("S14").Value = 0.08370346973784833
$ws.Range("T14").Value = 0.08370346973784833
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 74.41623166666666
$ws.Range("H15").Value = 223.248695
$ws.Range("I15").Value = 0.6076002364105371
$ws.Range("J15").Value = 0.6076002364105371
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 3.076282333333333
$ws.Range("N15").Value = 9.228847
$ws.Range("O15").Value = 0.2417686736584878
$ws.Range("P15").Value = 0.2417686736584878
$ws.Range("Q15").Value = 228.9253387894072
$ws.Range("R15").Value = 2060.328049104665
$ws.Range("S15").Value = 0.1468987032715592
$ws.Range("T15").Value = 0.1468987032715592
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 74.41623166666666
$ws.Range("H16").Value = 223.248695
$ws.Range("I16").Value = 0.6076002364105371
$ws.Range("J16").Value = 0.6076002364105371
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 7.303088666666667
$ws.Range("N16").Value = 21.909266
$ws.Range("O16").Value = 0.5739583917309499
$ws.Range("P16").Value = 0.5739583917309499
$ws.Range("Q16").Value = 543.4683381008743
$ws.Range("R16").Value = 4891.215042907869
$ws.Range("S16").Value = 0.3487372545055368
$ws.Range("T16").Value = 0.3487372545055368
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 74.41623166666666
$ws.Range("H17").Value = 223.248695
$ws.Range("I17").Value = 0.6076002364105371
$ws.Range("J17").Value = 0.6076002364105371
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.5918243333333334
$ws.Range("N17").Value = 1.775473
$ws.Range("O17").Value = 0.04651217560833507
$ws.Range("P17").Value = 0.04651217560833507
$ws.Range("Q17").Value = 44.04133669530389
$ws.Range("R17").Value = 396.372030257735
$ws.Range("S17").Value = 0.02826080889559281
$ws.Range("T17").Value = 0.02826080889559281
